# Apply the edit described by the diff:
# Insert 4 new data rows at position 271 (pushing existing rows 271-342 down to 275-346),
# and populate the new rows with the new observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before current row 271; this shifts rows 271:342 -> 275:346
$ws.Rows("271:274").Insert()

# Helper values that are constant across the whole data range
$marketId   = 8
$market     = "Terminal La Palmera de La Serena"
$region     = "Coquimbo"
$codreg     = 4
$tipo       = "Fruta"
$prodId     = 100103
$producto   = "Frutos de hueso (carozo)"
$catId      = 100103001
$categoria  = "Cereza"

# New row 271: Lapins / Especial
$r = 271
$ws.Cells.Item($r,1).Value  = $marketId
$ws.Cells.Item($r,2).Value  = $market
$ws.Cells.Item($r,3).Value  = $region
$ws.Cells.Item($r,4).Value  = 44900
$ws.Cells.Item($r,5).Value  = $codreg
$ws.Cells.Item($r,6).Value  = $tipo
$ws.Cells.Item($r,7).Value  = $prodId
$ws.Cells.Item($r,8).Value  = $producto
$ws.Cells.Item($r,9).Value  = $catId
$ws.Cells.Item($r,10).Value = $categoria
$ws.Cells.Item($r,11).Value = "Lapins"
$ws.Cells.Item($r,12).Value = "Especial"
$ws.Cells.Item($r,13).Value = 400
$ws.Cells.Item($r,14).Value = 9000
$ws.Cells.Item($r,15).Value = 10000
$ws.Cells.Item($r,16).Value = 9500
$ws.Cells.Item($r,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item($r,18).Value = "Región de O'Higgins"
$ws.Cells.Item($r,19).Value = 950
$ws.Cells.Item($r,20).Value = 10

# New row 272: Lapins / Primera
$r = 272
$ws.Cells.Item($r,1).Value  = $marketId
$ws.Cells.Item($r,2).Value  = $market
$ws.Cells.Item($r,3).Value  = $region
$ws.Cells.Item($r,4).Value  = 44900
$ws.Cells.Item($r,5).Value  = $codreg
$ws.Cells.Item($r,6).Value  = $tipo
$ws.Cells.Item($r,7).Value  = $prodId
$ws.Cells.Item($r,8).Value  = $producto
$ws.Cells.Item($r,9).Value  = $catId
$ws.Cells.Item($r,10).Value = $categoria
$ws.Cells.Item($r,11).Value = "Lapins"
$ws.Cells.Item($r,12).Value = "Primera"
$ws.Cells.Item($r,13).Value = 300
$ws.Cells.Item($r,14).Value = 7000
$ws.Cells.Item($r,15).Value = 8000
$ws.Cells.Item($r,16).Value = 7500
$ws.Cells.Item($r,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item($r,18).Value = "Región de O'Higgins"
$ws.Cells.Item($r,19).Value = 750
$ws.Cells.Item($r,20).Value = 10

# New row 273: Santina / Especial
$r = 273
$ws.Cells.Item($r,1).Value  = $marketId
$ws.Cells.Item($r,2).Value  = $market
$ws.Cells.Item($r,3).Value  = $region
$ws.Cells.Item($r,4).Value  = 44900
$ws.Cells.Item($r,5).Value  = $codreg
$ws.Cells.Item($r,6).Value  = $tipo
$ws.Cells.Item($r,7).Value  = $prodId
$ws.Cells.Item($r,8).Value  = $producto
$ws.Cells.Item($r,9).Value  = $catId
$ws.Cells.Item($r,10).Value = $categoria
$ws.Cells.Item($r,11).Value = "Santina"
$ws.Cells.Item($r,12).Value = "Especial"
$ws.Cells.Item($r,13).Value = 400
$ws.Cells.Item($r,14).Value = 10000
$ws.Cells.Item($r,15).Value = 11000
$ws.Cells.Item($r,16).Value = 10500
$ws.Cells.Item($r,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item($r,18).Value = "Provincia de Curicó"
$ws.Cells.Item($r,19).Value = 1050
$ws.Cells.Item($r,20).Value = 10

# New row 274: Santina / Primera
$r = 274
$ws.Cells.Item($r,1).Value  = $marketId
$ws.Cells.Item($r,2).Value  = $market
$ws.Cells.Item($r,3).Value  = $region
$ws.Cells.Item($r,4).Value  = 44900
$ws.Cells.Item($r,5).Value  = $codreg
$ws.Cells.Item($r,6).Value  = $tipo
$ws.Cells.Item($r,7).Value  = $prodId
$ws.Cells.Item($r,8).Value  = $producto
$ws.Cells.Item($r,9).Value  = $catId
$ws.Cells.Item($r,10).Value = $categoria
$ws.Cells.Item($r,11).Value = "Santina"
$ws.Cells.Item($r,12).Value = "Primera"
$ws.Cells.Item($r,13).Value = 300
$ws.Cells.Item($r,14).Value = 8000
$ws.Cells.Item($r,15).Value = 9000
$ws.Cells.Item($r,16).Value = 8500
$ws.Cells.Item($r,17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item($r,18).Value = "Provincia de Curicó"
$ws.Cells.Item($r,19).Value = 850
$ws.Cells.Item($r,20).Value = 10

# Make sure the date column keeps the date number format for the new rows
$ws.Range("D271:D274").NumberFormat = "YYYY-MM-DD HH:MM:SS"
